$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "tags"
$ws.Range("J2").Value = "testkey=testvalue;"

$ws.Range("J2").Select()
